# NSMB - 8-3 Done & some of 8-4
# Updates the speedrun split sheet ("V4") with the 8-3 finish time and
# the first part of the 8-4 splits.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("V4")
$ws.Activate()

# --- small corrections on existing rows ------------------------------
$ws.Range("B122").Value = 37222
$ws.Range("B123").Value = 37448

# extra observation values added alongside row 116
$ws.Range("J116").Value = 35610
$ws.Range("K116").Value = 35604

# --- row 124 gets filled in (previously only had C124) ---------------
$ws.Range("A124").Value = "Enter Pipe"
$ws.Range("B124").Value = 37664
# C124 already contains 44249

# --- brand new rows 125-132 (8-3 end .. start of 8-4) -----------------
$ws.Range("A125").Value = "Checkpoint Rail 12582912"
$ws.Range("B125").Value = 37919
$ws.Range("C125").Value = 44527

$ws.Range("A126").Value = "Checkpoint 460"
$ws.Range("B126").Value = 39165
$ws.Range("C126").Value = 45775

$ws.Range("A127").Value = "Cehckpiont 1449"
$ws.Range("B127").Value = 39722
$ws.Range("C127").Value = 46333

$ws.Range("A128").Value = "Enter Pipe"
$ws.Range("B128").Value = 41771
$ws.Range("C128").Value = 48401

$ws.Range("A129").Value = "Get flag"
$ws.Range("B129").Value = 42012
$ws.Range("C129").Value = 48651

$ws.Range("A130").Value = "End level"
$ws.Range("B130").Value = 42526
$ws.Range("C130").Value = 49165

$ws.Range("A131").Value = "Enter 8-4"
$ws.Range("B131").Value = 42919
$ws.Range("C131").Value = 49956

$ws.Range("A132").Value = "1st Move"
$ws.Range("B132").Value = 43147
$ws.Range("C132").Value = 50203

# --- extend the Diff column formula (IF(B>0,C-B,0)) down to row 132 --
$ws.Range("D124:D132").Formula = "=IF(B124>0,C124-B124,0)"

# --- refresh the frozen pane / selection to the new bottom of data ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 116
$ws.Range("B133").Select()
